$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.242.66"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.788.61"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'629.00"
$ws.Range("E5").Value = "  +4.59%  "
$ws.Range("D6").Value = "'164.62"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "3.787.79"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "'0.452"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'6.66"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'35.46"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "4.423.10"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "3.763.94"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "69.202.38"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "'17.96"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "'7.11"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'468.92"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'0.706"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'83.28"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'12.06"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "3.934.26"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "'2.24"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'7.32"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'29.03"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'9.03"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "3.737.09"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("E39").Value = "  +7.90%  "
$ws.Range("D40").Value = "'3.34"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.300"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'153.27"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'1.93"
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("D48").Value = "'46.89"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "'42.72"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'8.45"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("E51").Value = "  +2.99%  "
